$wb = $excel.ActiveWorkbook

# --- Sheet "Log": add new row 42 ---
$logSheet = $wb.Worksheets.Item("Log")
$logSheet.Range("A42").Value = "28/02/2025"
$logSheet.Range("B42").Value = "09:15"
$logSheet.Range("C42").Value = "Comisiones/Ventas % y modal By Categoría"
$logSheet.Range("D42").Value = "Ratio Comisiones/Ventas % incluye categoría Comisiones y Sueldos con descripción Comisiones Ventas (comision/comisones). Modal By Categoría ya usaba getCategoriaDisplay con la misma regla."
$logSheet.Range("E42").Value = "Diagnostico"

# --- Sheet "Versiones": add new row 13 ---
$versSheet = $wb.Worksheets.Item("Versiones")
$versSheet.Range("A13").NumberFormat = "@"
$versSheet.Range("A13").Value = "1.11"
$versSheet.Range("B13").Value = "28/02/2025"
$versSheet.Range("C13").Value = "Comisiones/Ventas %: categoría Comisiones + Sueldos (Comisiones Ventas); misma regla en modal By Categoría"
